$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 243, shifting rows 243:252 down to 244:253
$ws.Rows.Item(243).Insert()

# Fill the new row 243 with its data values (columns B..J)
$ws.Cells.Item(243, 2).Value = "1A"
$ws.Cells.Item(243, 3).Value = "1B"
$ws.Cells.Item(243, 4).Value = "2A"
$ws.Cells.Item(243, 5).Value = "2B"
$ws.Cells.Item(243, 6).Value = "2A"
$ws.Cells.Item(243, 7).Value = "2B"

# Columns H and I hold text category codes ("3" and "4") that look like plain
# numbers. Assigning ".Value = "3"" would make Excel store them as numeric
# values instead of text, so instead copy them from existing cells in the
# sheet that already hold those same values as text (H2 = "3", I2 = "4").
$ws.Cells.Item(2, 8).Copy($ws.Cells.Item(243, 8))
$ws.Cells.Item(2, 9).Copy($ws.Cells.Item(243, 9))

$ws.Cells.Item(243, 10).Value = "1A"

# Renumber column A (sequential record numbers) for rows 2..253
for ($r = 2; $r -le 253; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}
